# Update crypto price/volume data rows per Fri Apr 21 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.593.57'
$ws.Cells.Item(2, 5).Value = '  -1.56%  '

$ws.Cells.Item(3, 4).Value = '1.967.36'
$ws.Cells.Item(3, 5).Value = '  +0.32%  '

$ws.Cells.Item(4, 5).Value = '  +0.29%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '323.93'
$ws.Cells.Item(5, 5).Value = '  -0.94%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.011'
$ws.Cells.Item(6, 5).Value = '  +0.36%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4805'
$ws.Cells.Item(7, 5).Value = '  -3.51%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4073'
$ws.Cells.Item(8, 5).Value = '  -3.24%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '54.03'
$ws.Cells.Item(9, 5).Value = '  -0.21%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08528'
$ws.Cells.Item(10, 5).Value = '  -5.14%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.065'
$ws.Cells.Item(11, 5).Value = '  -2.97%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.49'
$ws.Cells.Item(12, 5).Value = '  -2.17%  '

$ws.Cells.Item(13, 4).Value = '1.976.64'
$ws.Cells.Item(13, 5).Value = '  +1.32%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.644'
$ws.Cells.Item(14, 5).Value = '  -2.67%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.209'
$ws.Cells.Item(15, 5).Value = '  -3.57%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.013'
$ws.Cells.Item(16, 5).Value = '  +0.41%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '91.41'
$ws.Cells.Item(17, 5).Value = '  +0.33%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001076'
$ws.Cells.Item(18, 5).Value = '  -1.79%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06658'
$ws.Cells.Item(19, 5).Value = '  -0.10%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.66'
$ws.Cells.Item(20, 5).Value = '  -2.54%  '

$ws.Cells.Item(21, 5).Value = '  +0.40%  '

$ws.Cells.Item(22, 5).Value = '  -1.11%  '

$ws.Cells.Item(23, 4).Value = '28.613.29'
$ws.Cells.Item(23, 5).Value = '  -1.53%  '

$ws.Cells.Item(24, 5).Value = '  -3.19%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.302'
$ws.Cells.Item(25, 5).Value = '  +0.66%  '

$ws.Cells.Item(26, 4).Value = '2.241.94'
$ws.Cells.Item(26, 5).Value = '  +2.35%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '156.49'
$ws.Cells.Item(27, 5).Value = '  +0.53%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '20.41'
$ws.Cells.Item(28, 5).Value = '  -0.82%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.908'
$ws.Cells.Item(29, 5).Value = '  -4.30%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.186'
$ws.Cells.Item(30, 5).Value = '  -2.67%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '125.02'
$ws.Cells.Item(31, 5).Value = '  -1.71%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.9931'
$ws.Cells.Item(32, 5).Value = '  -4.49%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.09690'
$ws.Cells.Item(33, 5).Value = '  -1.36%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.463'
$ws.Cells.Item(34, 5).Value = '  -4.31%  '

$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.707'
$ws.Cells.Item(35, 5).Value = '  +0.32%  '

$ws.Cells.Item(36, 2).Value = 'Filecoin'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.648'
$ws.Cells.Item(36, 5).Value = '  -2.42%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '9.166'
$ws.Cells.Item(37, 5).Value = '  +2.49%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02342'
$ws.Cells.Item(38, 5).Value = '  -2.94%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.06251'
$ws.Cells.Item(39, 5).Value = '  -0.54%  '

$ws.Cells.Item(40, 5).Value = '  -2.31%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6251'
$ws.Cells.Item(41, 5).Value = '  -2.63%  '

$ws.Cells.Item(42, 5).Value = '  -1.59%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.011'
$ws.Cells.Item(43, 5).Value = '  +0.32%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1925'
$ws.Cells.Item(44, 5).Value = '  -2.76%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.365'
$ws.Cells.Item(45, 5).Value = '  +6.77%  '

$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '13.16'
$ws.Cells.Item(46, 5).Value = '  -2.32%  '

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5978'
$ws.Cells.Item(47, 5).Value = '  -3.16%  '

$ws.Cells.Item(48, 5).Value = '  -4.39%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '3.417'
$ws.Cells.Item(49, 5).Value = '  -1.37%  '

$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.06834'
$ws.Cells.Item(50, 5).Value = '  -0.59%  '

$ws.Cells.Item(51, 2).Value = 'Quant'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '111.67'
$ws.Cells.Item(51, 5).Value = '  -0.64%  '

